$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------------
# 1) Rewrite the "Alerts" user-story paragraph text (row labelled "Alerts").
#    Do it as three scoped Find/Replace passes on unique anchor phrases so
#    the similarly-worded paragraph in row 15 ("...alerts for delayed and
#    cancelled buses...") is left untouched.
# ---------------------------------------------------------------------------

# "specific alerts for next arriving buses and  general alerts for" -> "alerts for the buses"
$found1 = $d.Content.Find.Execute("specific alerts for next arriving buses and  general alerts for", $true, $false, $false, $false, $false, $true, 1, $false, "alerts for the buses", 2)

# "delayed/cancelled buses so " -> "so " (drop the "delayed/cancelled buses " run)
$found2 = $d.Content.Find.Execute("delayed/cancelled buses so ", $true, $false, $false, $false, $false, $true, 1, $false, "so ", 2)

# "that I can spend less time waiting for the bus and change my plan
#  according to the real-time route status" -> "that know the real-time
#  status of the buses."
$search3 = "that I can" + $nbsp + "spend less time waiting for the bus and" + $nbsp + "change my plan according to the real-time route status"
$found3 = $d.Content.Find.Execute($search3, $true, $false, $false, $false, $false, $true, 1, $false, "that know the real-time status of the buses.", 2)

# ---------------------------------------------------------------------------
# 2) Insert a new (collapsed) "_GoBack" bookmark right before "buses" in the
#    rewritten sentence ("...for the |buses...").
# ---------------------------------------------------------------------------
$r = $d.Content
$found4 = $r.Find.Execute("for the buses", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $r.Start + 8
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3) Remove the stale "_GoBack" bookmark that used to sit between
#    "...bus timings " and "so that..." in the next row. "_GoBack" is a
#    hidden/system bookmark so it cannot be reached through the Bookmarks
#    collection; re-typing the text that spans it drops it as a side effect.
# ---------------------------------------------------------------------------
$found5 = $d.Content.Find.Execute("bus timings so that", $true, $false, $false, $false, $false, $true, 1, $false, "bus timings so that", 2)

# ---------------------------------------------------------------------------
# 4) Drop the stale <w:lastRenderedPageBreak/> rendering hint in front of the
#    "16" row number. Re-typing the cell's text removes the hint element
#    while preserving its run formatting (yellow highlight).
# ---------------------------------------------------------------------------
$found6 = $d.Content.Find.Execute("16", $true, $true, $false, $false, $false, $true, 1, $false, "16", 2)

Write-Host "replace1:" $found1 " replace2:" $found2 " replace3:" $found3 " bookmarkFind:" $found4 " bookmarkCleanup:" $found5 " pageBreakCleanup:" $found6
